# Auto-generated from XML diff: update crypto price table rows 2-51
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = '20.009.90'
$ws.Range("E2").Value = '  -4.60%  '

# Row 3
$ws.Range("D3").Value = '1.419.81'
$ws.Range("E3").Value = '  -5.07%  '

# Row 4
$ws.Range("D4").Value = "'0.9981"
$ws.Range("E4").Value = '  -0.93%  '

# Row 5
$ws.Range("D5").Value = "'0.9997"
$ws.Range("E5").Value = '  -0.64%  '

# Row 6
$ws.Range("D6").Value = "'276.74"
$ws.Range("E6").Value = '  -1.93%  '

# Row 7
$ws.Range("D7").Value = "'0.3671"
$ws.Range("E7").Value = '  -3.48%  '

# Row 8
$ws.Range("D8").Value = "'0.3102"
$ws.Range("E8").Value = '  -0.37%  '

# Row 9
$ws.Range("D9").Value = "'39.75"
$ws.Range("E9").Value = '  -7.01%  '

# Row 10
$ws.Range("D10").Value = "'1.045"
$ws.Range("E10").Value = '  +1.99%  '

# Row 11
$ws.Range("D11").Value = "'0.06555"
$ws.Range("E11").Value = '  -4.69%  '

# Row 12
$ws.Range("D12").Value = "'0.9979"
$ws.Range("E12").Value = '  -1.02%  '

# Row 13
$ws.Range("D13").Value = "'5.513"
$ws.Range("E13").Value = '  -0.24%  '

# Row 14
$ws.Range("D14").Value = "'17.75"
$ws.Range("E14").Value = '  +0.54%  '

# Row 15
$ws.Range("D15").Value = "'6.213"
$ws.Range("E15").Value = '  -1.80%  '

# Row 16
$ws.Range("B16").Value = 'ShibaInu'
$ws.Range("C16").Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range("D16").Value = "'0.00001023"
$ws.Range("E16").Value = '  -3.40%  '

# Row 17
$ws.Range("B17").Value = 'WrappedEther'
$ws.Range("C17").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D17").Value = '1.416.21'
$ws.Range("E17").Value = '  -5.93%  '

# Row 18
$ws.Range("D18").Value = "'0.05677"
$ws.Range("E18").Value = '  -13.25%  '

# Row 19
$ws.Range("D19").Value = "'0.9993"
$ws.Range("E19").Value = '  -0.66%  '

# Row 20
$ws.Range("D20").Value = "'71.39"
$ws.Range("E20").Value = '  -12.54%  '

# Row 21
$ws.Range("D21").Value = "'5.634"
$ws.Range("E21").Value = '  -4.99%  '

# Row 22
$ws.Range("D22").Value = "'14.78"
$ws.Range("E22").Value = '  -1.18%  '

# Row 23
$ws.Range("D23").Value = "'11.02"
$ws.Range("E23").Value = '  +2.05%  '

# Row 24
$ws.Range("D24").Value = "'2.249"
$ws.Range("E24").Value = '  -3.79%  '

# Row 25
$ws.Range("D25").Value = '20.003.62'
$ws.Range("E25").Value = '  -4.68%  '

# Row 26
$ws.Range("D26").Value = "'2.283"
$ws.Range("E26").Value = '  -1.03%  '

# Row 27
$ws.Range("D27").Value = "'133.40"
$ws.Range("E27").Value = '  -8.99%  '

# Row 28
$ws.Range("D28").Value = "'17.40"
$ws.Range("E28").Value = '  -2.70%  '

# Row 29
$ws.Range("D29").Value = '1.576.44'
$ws.Range("E29").Value = '  -5.82%  '

# Row 30
$ws.Range("D30").Value = "'110.03"
$ws.Range("E30").Value = '  -3.20%  '

# Row 31
$ws.Range("D31").Value = "'3.894"
$ws.Range("E31").Value = '  -18.64%  '

# Row 32
$ws.Range("D32").Value = "'5.290"
$ws.Range("E32").Value = '  -8.56%  '

# Row 33
$ws.Range("D33").Value = "'0.8253"
$ws.Range("E33").Value = '  -11.97%  '

# Row 34
$ws.Range("D34").Value = "'0.07713"
$ws.Range("E34").Value = '  -2.17%  '

# Row 35
$ws.Range("E35").Value = '  +0.91%  '

# Row 36
$ws.Range("D36").Value = "'8.311"
$ws.Range("E36").Value = '  -0.55%  '

# Row 37
$ws.Range("D37").Value = "'4.950"
$ws.Range("E37").Value = '  -1.00%  '

# Row 38
$ws.Range("D38").Value = "'0.05825"
$ws.Range("E38").Value = '  +1.70%  '

# Row 39
$ws.Range("D39").Value = "'0.9994"
$ws.Range("E39").Value = '  -0.60%  '

# Row 40
$ws.Range("D40").Value = "'0.02065"
$ws.Range("E40").Value = '  -1.95%  '

# Row 41
$ws.Range("D41").Value = "'10.55"
$ws.Range("E41").Value = '  -3.27%  '

# Row 42
$ws.Range("D42").Value = "'0.1886"
$ws.Range("E42").Value = '  -3.52%  '

# Row 43
$ws.Range("D43").Value = "'1.102"
$ws.Range("E43").Value = '  -4.45%  '

# Row 44
$ws.Range("B44").Value = 'TheSandbox'
$ws.Range("C44").Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range("D44").Value = "'0.5338"
$ws.Range("E44").Value = '  -3.97%  '

# Row 45
$ws.Range("B45").Value = 'EnergySwap'
$ws.Range("C45").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D45").Value = "'12.39"
$ws.Range("E45").Value = '  -2.49%  '

# Row 46
$ws.Range("D46").Value = "'3.543"
$ws.Range("E46").Value = '  -2.77%  '

# Row 47
$ws.Range("D47").Value = "'0.5207"
$ws.Range("E47").Value = '  -3.31%  '

# Row 48
$ws.Range("D48").Value = "'116.37"
$ws.Range("E48").Value = '  +3.45%  '

# Row 49
$ws.Range("D49").Value = "'1.780"
$ws.Range("E49").Value = '  -2.35%  '

# Row 50
$ws.Range("D50").Value = "'1.035"
$ws.Range("E50").Value = '  -7.66%  '

# Row 51
$ws.Range("E51").Value = '  -0.60%  '
